$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-28 03:18:13'
$ws.Range('N2').Value = '0.4 °C 2:51 TU'
$ws.Range('O2').Value = '1.0 °C'
$ws.Range('E3').Value = '2026-02-28 03:18:15'
$ws.Range('H3').NumberFormat = "@"
$ws.Range('H3').Value = '79%'
$ws.Range('O3').Value = '0.0 °C'
$ws.Range('E4').Value = '2026-02-28 03:18:17'
$ws.Range('J4').Value = '1024.4 hPa'
$ws.Range('L4').Value = '7.6 km/h - 307º 2:48 TU'
$ws.Range('E5').Value = '2026-02-28 03:18:20'
$ws.Range('H5').NumberFormat = "@"
$ws.Range('H5').Value = '98%'
$ws.Range('L5').Value = '15.1 km/h - 55º 2:52 TU'
$ws.Range('N5').Value = '-0.3 °C 2:51 TU'
$ws.Range('O5').Value = '0.1 °C'
$ws.Range('E6').Value = '2026-02-28 03:18:22'
$ws.Range('E7').Value = '2026-02-28 03:18:24'
$ws.Range('J7').Value = '1023.8 hPa'
$ws.Range('L7').Value = '18.7 km/h - 48º 2:54 TU'
$ws.Range('E8').Value = '2026-02-28 03:18:27'
$ws.Range('L8').Value = '13.0 km/h - 47º 2:58 TU'
$ws.Range('M8').Value = '8.6 °C 2:43 TU'
$ws.Range('O8').Value = '8.4 °C'
$ws.Range('E9').Value = '2026-02-28 03:18:29'
$ws.Range('L9').Value = '5.4 km/h - 340º 2:59 TU'
$ws.Range('M9').Value = '8.2 °C 2:59 TU'
$ws.Range('O9').Value = '7.1 °C'
$ws.Range('E10').Value = '2026-02-28 03:18:30'
$ws.Range('M10').Value = '8.3 °C 2:53 TU'
$ws.Range('O10').Value = '7.3 °C'
$ws.Range('E11').Value = '2026-02-28 03:18:31'
$ws.Range('H11').NumberFormat = "@"
$ws.Range('H11').Value = '94%'
$ws.Range('N11').Value = '3.3 °C 2:56 TU'
$ws.Range('O11').Value = '4.0 °C'
$ws.Range('E12').Value = '2026-02-28 03:18:32'
$ws.Range('E13').Value = '2026-02-28 03:18:35'
$ws.Range('G13').Value = '3 cm'
$ws.Range('H13').NumberFormat = "@"
$ws.Range('H13').Value = '83%'
$ws.Range('J13').Value = '1025.7 hPa'
$ws.Range('N13').Value = '0.6 °C 2:40 TU'
$ws.Range('O13').Value = '2.2 °C'
$ws.Range('E14').Value = '2026-02-28 03:18:37'
$ws.Range('O14').Value = '9.8 °C'
$ws.Range('E15').Value = '2026-02-28 03:18:39'
$ws.Range('M15').Value = '7.7 °C 2:52 TU'
$ws.Range('O15').Value = '6.4 °C'
$ws.Range('E16').Value = '2026-02-28 03:18:42'
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H16').Value = '65%'
$ws.Range('N16').Value = '-1.0 °C 2:59 TU'
$ws.Range('O16').Value = '-0.4 °C'
$ws.Range('E17').Value = '2026-02-28 03:18:44'
$ws.Range('H17').NumberFormat = "@"
$ws.Range('H17').Value = '45%'
$ws.Range('E18').Value = '2026-02-28 03:18:46'
$ws.Range('J18').Value = '1024.4 hPa'
$ws.Range('M18').Value = '8.4 °C 2:59 TU'
$ws.Range('O18').Value = '7.7 °C'
$ws.Range('E19').Value = '2026-02-28 03:18:49'
$ws.Range('E20').Value = '2026-02-28 03:18:51'
$ws.Range('H20').NumberFormat = "@"
$ws.Range('H20').Value = '37%'
$ws.Range('L20').Value = '19.1 km/h - 173º 2:59 TU'
$ws.Range('N20').Value = '-0.5 °C 2:50 TU'
$ws.Range('O20').Value = '0.4 °C'
$ws.Range('E21').Value = '2026-02-28 03:18:53'
$ws.Range('N21').Value = '5.1 °C 2:57 TU'
$ws.Range('E22').Value = '2026-02-28 03:18:56'
$ws.Range('H22').NumberFormat = "@"
$ws.Range('H22').Value = '62%'
$ws.Range('O22').Value = '-0.5 °C'
$ws.Range('E23').Value = '2026-02-28 03:18:58'
$ws.Range('H23').NumberFormat = "@"
$ws.Range('H23').Value = '67%'
$ws.Range('N23').Value = '-0.6 °C 2:34 TU'
$ws.Range('O23').Value = '0.4 °C'
$ws.Range('E24').Value = '2026-02-28 03:19:00'
$ws.Range('H24').NumberFormat = "@"
$ws.Range('H24').Value = '97%'
$ws.Range('J24').Value = '1023.2 hPa'
$ws.Range('N24').Value = '4.4 °C 2:59 TU'
$ws.Range('O24').Value = '6.8 °C'
$ws.Range('E25').Value = '2026-02-28 03:19:03'
$ws.Range('H25').NumberFormat = "@"
$ws.Range('H25').Value = '58%'
$ws.Range('N25').Value = '0.2 °C 2:38 TU'
$ws.Range('O25').Value = '1.5 °C'
$ws.Range('E26').Value = '2026-02-28 03:19:05'
$ws.Range('H26').NumberFormat = "@"
$ws.Range('H26').Value = '71%'
$ws.Range('E27').Value = '2026-02-28 03:19:07'
$ws.Range('N27').Value = '1.9 °C 2:59 TU'
$ws.Range('O27').Value = '3.0 °C'
$ws.Range('E28').Value = '2026-02-28 03:19:10'
$ws.Range('H28').NumberFormat = "@"
$ws.Range('H28').Value = '96%'
$ws.Range('J28').Value = '1024.5 hPa'
$ws.Range('N28').Value = '5.5 °C 2:58 TU'
$ws.Range('O28').Value = '6.7 °C'
$ws.Range('E29').Value = '2026-02-28 03:19:12'
$ws.Range('O29').Value = '8.3 °C'
$ws.Range('E30').Value = '2026-02-28 03:19:14'
$ws.Range('E31').Value = '2026-02-28 03:19:17'
$ws.Range('H31').NumberFormat = "@"
$ws.Range('H31').Value = '93%'
$ws.Range('L31').Value = '36.7 km/h - 10º 2:58 TU'
$ws.Range('N31').Value = '9.9 °C 2:56 TU'
$ws.Range('O31').Value = '10.3 °C'
$ws.Range('E32').Value = '2026-02-28 03:19:19'
$ws.Range('H32').NumberFormat = "@"
$ws.Range('H32').Value = '83%'
$ws.Range('O32').Value = '5.0 °C'
$ws.Range('E33').Value = '2026-02-28 03:19:22'
$ws.Range('J33').Value = '1023.4 hPa'
$ws.Range('N33').Value = '4.1 °C 2:57 TU'
$ws.Range('O33').Value = '5.5 °C'
$ws.Range('E34').Value = '2026-02-28 03:19:24'
$ws.Range('H34').NumberFormat = "@"
$ws.Range('H34').Value = '73%'
$ws.Range('N34').Value = '-0.6 °C 2:30 TU'
$ws.Range('O34').Value = '0.2 °C'
$ws.Range('E35').Value = '2026-02-28 03:19:26'
$ws.Range('L35').Value = '25.9 km/h - 239º 2:58 TU'
$ws.Range('O35').Value = '7.2 °C'
$ws.Range('E36').Value = '2026-02-28 03:19:28'
$ws.Range('L36').Value = '7.6 km/h - 252º 2:37 TU'
$ws.Range('E37').Value = '2026-02-28 03:19:31'
$ws.Range('H37').NumberFormat = "@"
$ws.Range('H37').Value = '87%'
$ws.Range('L37').Value = '12.6 km/h - 39º 2:36 TU'
$ws.Range('E38').Value = '2026-02-28 03:19:33'
$ws.Range('H38').NumberFormat = "@"
$ws.Range('H38').Value = '91%'
$ws.Range('O38').Value = '8.8 °C'
$ws.Range('E39').Value = '2026-02-28 03:19:35'
$ws.Range('H39').NumberFormat = "@"
$ws.Range('H39').Value = '52%'
$ws.Range('O39').Value = '0.7 °C'
$ws.Range('E40').Value = '2026-02-28 03:19:37'
$ws.Range('H40').NumberFormat = "@"
$ws.Range('H40').Value = '92%'
$ws.Range('O40').Value = '3.8 °C'
$ws.Range('E41').Value = '2026-02-28 03:19:40'
$ws.Range('N41').Value = '11.7 °C 2:30 TU'
$ws.Range('E42').Value = '2026-02-28 03:19:42'
$ws.Range('M42').Value = '8.0 °C 2:57 TU'
$ws.Range('O42').Value = '7.5 °C'
$ws.Range('E43').Value = '2026-02-28 03:19:44'
$ws.Range('H43').NumberFormat = "@"
$ws.Range('H43').Value = '82%'
$ws.Range('O43').Value = '4.1 °C'
$ws.Range('E44').Value = '2026-02-28 03:19:47'
$ws.Range('H44').NumberFormat = "@"
$ws.Range('H44').Value = '89%'
$ws.Range('O44').Value = '-0.8 °C'
$ws.Range('E45').Value = '2026-02-28 03:19:49'
$ws.Range('H45').NumberFormat = "@"
$ws.Range('H45').Value = '86%'
$ws.Range('N45').Value = '6.7 °C 2:58 TU'
$ws.Range('O45').Value = '7.6 °C'
$ws.Range('E46').Value = '2026-02-28 03:19:51'
$ws.Range('O46').Value = '10.9 °C'
